$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "run" at G1
$ws.Range("G1").Value = "run"

# Row 4 (Loïse): add F4 and G4 = "X"
$ws.Range("F4").Value = "X"
$ws.Range("G4").Value = "X"

# Rows 5-8: add G column = "X"
$ws.Range("G5").Value = "X"
$ws.Range("G6").Value = "X"
$ws.Range("G7").Value = "X"
$ws.Range("G8").Value = "X"

# Update selection to G4 as in the diff
$ws.Range("G4").Select()
